$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - add Resume_received / Resume_downloaded values
$ws.Range("H2").Value = 12
$ws.Range("I2").Value = 12

# Row 3 - add Resume_received / Resume_downloaded values
$ws.Range("H3").Value = 23
$ws.Range("I3").Value = 23

# Row 4 - add LinkedIn_Poster / LinkedIn_Posted / Resume_received / Resume_downloaded values
$ws.Range("F4").Value = "Created"
$ws.Range("G4").Value = "Yes"
$ws.Range("H4").Value = 44
$ws.Range("I4").Value = 44

# Row 5 - add LinkedIn_Poster / LinkedIn_Posted / Resume_received / Resume_downloaded values
$ws.Range("F5").Value = "Created"
$ws.Range("G5").Value = "Yes"
$ws.Range("H5").Value = 11
$ws.Range("I5").Value = 11

# Update the active selection to J6
$ws.Range("J6").Select()
